# Regenerate orders with updated distance/size codes.
# Rule (derived from the diff): within the text of the relevant columns,
# replace the old distance/size tokens with the new ones.
#   D51 -> D55
#   D80 -> D86
#   D64 -> D69
#   S30 -> S31
# (S25/S20 are unchanged.)  Applied to the Condition, Filename_Left,
# Filename_Right, Distance and Size columns - the only columns whose
# shared-string text actually contains these tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Header row is row 1; data starts at row 2.
$targetCols = @(2, 4, 5, 8, 10)  # Condition, Filename_Left, Filename_Right, Distance, Size

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $newVal = $val.Replace("D51", "D55").Replace("D80", "D86").Replace("D64", "D69").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
